$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.779.41'
$ws.Range('E2').Value = '  +7.07%  '
$ws.Range('D3').Value = '3.329.37'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '410.51'
$ws.Range('E5').Value = '  +3.71%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '117.43'
$ws.Range('E6').Value = '  +7.65%  '
$ws.Range('D7').Value = '3.322.98'
$ws.Range('E7').Value = '  +2.31%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.574'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.629'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').Value = '  +15.94%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '40.30'
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').Value = '3.848.20'
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.24'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.16'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '3.324.60'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').Value = '60.557.66'
$ws.Range('E18').Value = '  +6.85%  '
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.89'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.38'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('E22').Value = '  +4.90%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.42'
$ws.Range('E23').Value = '  -4.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '297.03'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.02'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '29.07'
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.82'
$ws.Range('E28').Value = '  +6.59%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.22'
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.171'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.49'
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('E32').Value = '  +4.45%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '42.53'
$ws.Range('E33').Value = '  +7.32%  '
$ws.Range('B34').Value = 'Dai'
$ws.Range('C34').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('B35').Value = 'Cosmos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '11.29'
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.51'
$ws.Range('E36').Value = '  +17.58%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0491'
$ws.Range('E37').Value = '  +0.85%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '52.13'
$ws.Range('E38').Value = '  +1.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.997'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.07'
$ws.Range('E40').Value = '  +5.93%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '133.44'
$ws.Range('E42').Value = '  -4.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.291'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.121'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.90'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.93'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.31'
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('E48').Value = '  +4.13%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '21.20'
$ws.Range('E49').Value = '  -5.01%  '
$ws.Range('D50').Value = '2.145.63'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('D51').Value = '3.653.85'
$ws.Range('E51').Value = '  +2.27%  '
